$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-65 down to 25-66.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly price entry.
$ws.Range("A24").Value = 1
$ws.Range("B24").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C24").Value = "Arica y Parinacota"
$ws.Range("D24").Value = 45274
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = 100112044
$ws.Range("G24").Value = "Perejil"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 450
$ws.Range("K24").Value = 800
$ws.Range("L24").Value = 1000
$ws.Range("M24").Value = 911
$ws.Range("N24").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O24").Value = "Región de Arica y Parinacota"
$ws.Range("P24").Value = 456
$ws.Range("Q24").Value = 2
$ws.Range("R24").Value = "Hortaliza"

Write-Host "Done"
